$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Paragraph 1 (the question): remove the _GoBack bookmark and
# split "Q: Can I apply to CUT while completing A level?" into
# three runs: "Q: ... completing " / "A-level" / "?"
# ---------------------------------------------------------------

# The _GoBack bookmark currently wraps the very start of paragraph 1.
# Delete it outright (it gets re-created later, inside paragraph 2).
$d.Bookmarks("_GoBack").Delete()

$p1 = $d.Paragraphs(1).Range
$p1start = $p1.Start
$p1end = $p1.End
$q1full = $d.Range($p1start, $p1end)
$q1full.Text = "Q: Can I apply to CUT while completing "

$p1b = $d.Paragraphs(1).Range
$ins1 = $d.Range($p1b.End - 1, $p1b.End - 1)
$ins1.InsertAfter("A-level")

$p1c = $d.Paragraphs(1).Range
$ins2 = $d.Range($p1c.End - 1, $p1c.End - 1)
$ins2.InsertAfter("?")

# ---------------------------------------------------------------
# Paragraph 2 (the answer): replace the whole answer text and
# re-insert the _GoBack bookmark in the middle of the new text,
# between "...avoid regrets and inconve" and "nience, ...".
# ---------------------------------------------------------------

$p2 = $d.Paragraphs(2).Range
$p2start = $p2.Start
$p2end = $p2.End
$a1full = $d.Range($p2start, $p2end)
$a1full.Text = "A: The pre-application indicates that the prospective candidate is enthusiastic about our institution, which is greatly appreciated. However, admission into the desired program depends on submitting the results after successfully completing the qualifying studies. When applying, the system asks you to submit your qualification, and if it is not sent in full, you will be rejected. Therefore, to avoid regrets and inconve"

$p2b = $d.Paragraphs(2).Range
$bmPos = $d.Range($p2b.End - 1, $p2b.End - 1)
$d.Bookmarks.Add("_GoBack", $bmPos)

$p2c = $d.Paragraphs(2).Range
$ins3 = $d.Range($p2c.End - 1, $p2c.End - 1)
$ins3.InsertAfter("nience, pre-application is generally not recommended.")

Write-Host "Q:" $d.Paragraphs(1).Range.Text
Write-Host "A:" $d.Paragraphs(2).Range.Text
